# Regenerate save_data G column ("K") values: recomputed strike-count (K)
# values replacing the old "Strike#" derived figures, per regen of std/mean
# and s_vals calculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 2
    4  = 0
    5  = 1
    6  = 1
    7  = 2
    8  = 0
    9  = 1
    10 = 2
    11 = 0
    12 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 1
    19 = 0
    20 = 0
    21 = 1
    22 = 1
    23 = 0
    24 = 0
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 0
    30 = 0
    31 = 0
    32 = 1
    33 = 1
    34 = 1
    35 = 1
    36 = 3
    37 = 1
    38 = 3
    39 = 2
    40 = 0
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
